# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 (header) ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 15:53:28'
$ws1.Range("A3").Value = 'Total filas: 328'

# --- Sheet1 data block: A6:E333 (328 rows) ---
$data1 = New-Object 'object[,]' 328,5
$data1[0,0] = '03:45:25'
$data1[0,1] = '03:47'
$data1[0,2] = '14_ABASTO'
$data1[0,3] = 2
$data1[0,4] = 'LP1912'
$data1[1,0] = '03:45:25'
$data1[1,1] = '04:01'
$data1[1,2] = '81_EL PELIGRO'
$data1[1,3] = 16
$data1[1,4] = 'LP1912'
$data1[2,0] = '03:45:25'
$data1[2,1] = '04:46'
$data1[2,2] = '215A_EL PATO'
$data1[2,3] = 61
$data1[2,4] = 'LP1912'
$data1[3,0] = '03:45:25'
$data1[3,1] = '04:53'
$data1[3,2] = '11_ETCHEVERRY'
$data1[3,3] = 68
$data1[3,4] = 'LP1912'
$data1[4,0] = '04:56:49'
$data1[4,1] = '05:13'
$data1[4,2] = '14_ABASTO'
$data1[4,3] = 17
$data1[4,4] = 'LP1912'
$data1[5,0] = '03:45:25'
$data1[5,1] = '05:16'
$data1[5,2] = '17_ROMERO'
$data1[5,3] = 91
$data1[5,4] = 'LP1912'
$data1[6,0] = '04:45:05'
$data1[6,1] = '05:16'
$data1[6,2] = '14_ABASTO'
$data1[6,3] = 31
$data1[6,4] = 'LP1912'
$data1[7,0] = '03:45:25'
$data1[7,1] = '05:22'
$data1[7,2] = '23_HERNANDEZ'
$data1[7,3] = 97
$data1[7,4] = 'LP1912'
$data1[8,0] = '05:26:08'
$data1[8,1] = '05:28'
$data1[8,2] = '14_ABASTO'
$data1[8,3] = 2
$data1[8,4] = 'LP1912'
$data1[9,0] = '04:18:02'
$data1[9,1] = '05:34'
$data1[9,2] = '14_ABASTO'
$data1[9,3] = 76
$data1[9,4] = 'LP1912'
$data1[10,0] = '03:45:25'
$data1[10,1] = '05:34'
$data1[10,2] = '215B_EL PATO'
$data1[10,3] = 109
$data1[10,4] = 'LP1912'
$data1[11,0] = '04:18:02'
$data1[11,1] = '05:35'
$data1[11,2] = '215B_EL PATO'
$data1[11,3] = 77
$data1[11,4] = 'LP1912'
$data1[12,0] = '03:45:25'
$data1[12,1] = '05:37'
$data1[12,2] = '14_ABASTO'
$data1[12,3] = 112
$data1[12,4] = 'LP1912'
$data1[13,0] = '04:18:02'
$data1[13,1] = '05:46'
$data1[13,2] = '15_ABASTO'
$data1[13,3] = 88
$data1[13,4] = 'LP1912'
$data1[14,0] = '04:45:05'
$data1[14,1] = '06:04'
$data1[14,2] = '16_SANTA ANA'
$data1[14,3] = 79
$data1[14,4] = 'LP1912'
$data1[15,0] = '04:18:02'
$data1[15,1] = '06:05'
$data1[15,2] = '16_SANTA ANA'
$data1[15,3] = 107
$data1[15,4] = 'LP1912'
$data1[16,0] = '04:56:49'
$data1[16,1] = '06:11'
$data1[16,2] = '215A_EL PATO'
$data1[16,3] = 75
$data1[16,4] = 'LP1912'
$data1[17,0] = '04:18:02'
$data1[17,1] = '06:12'
$data1[17,2] = '215A_EL PATO'
$data1[17,3] = 114
$data1[17,4] = 'LP1912'
$data1[18,0] = '04:18:02'
$data1[18,1] = '06:14'
$data1[18,2] = '225_HARAS DEL SUR'
$data1[18,3] = 116
$data1[18,4] = 'LP1912'
$data1[19,0] = '04:45:05'
$data1[19,1] = '06:21'
$data1[19,2] = '26_HERNANDEZ'
$data1[19,3] = 96
$data1[19,4] = 'LP1912'
$data1[20,0] = '06:25:43'
$data1[20,1] = '06:26'
$data1[20,2] = '86_EST CHICA-ESC AGRARIA'
$data1[20,3] = 1
$data1[20,4] = 'LP1912'
$data1[21,0] = '04:45:05'
$data1[21,1] = '06:27'
$data1[21,2] = '23_HERNANDEZ'
$data1[21,3] = 102
$data1[21,4] = 'LP1912'
$data1[22,0] = '06:25:43'
$data1[22,1] = '06:28'
$data1[22,2] = '23_HERNANDEZ'
$data1[22,3] = 3
$data1[22,4] = 'LP1912'
$data1[23,0] = '04:56:49'
$data1[23,1] = '06:29'
$data1[23,2] = '86_EST CHICA-ESC AGRARIA'
$data1[23,3] = 93
$data1[23,4] = 'LP1912'
$data1[24,0] = '04:45:05'
$data1[24,1] = '06:30'
$data1[24,2] = '86_EST CHICA-ESC AGRARIA'
$data1[24,3] = 105
$data1[24,4] = 'LP1912'
$data1[25,0] = '04:45:05'
$data1[25,1] = '06:31'
$data1[25,2] = '16_SANTA ANA'
$data1[25,3] = 106
$data1[25,4] = 'LP1912'
$data1[26,0] = '05:55:25'
$data1[26,1] = '06:44'
$data1[26,2] = '26_HERNANDEZ'
$data1[26,3] = 49
$data1[26,4] = 'LP1912'
$data1[27,0] = '04:45:05'
$data1[27,1] = '06:44'
$data1[27,2] = '225_C ROCA-H SUR'
$data1[27,3] = 119
$data1[27,4] = 'LP1912'
$data1[28,0] = '04:56:49'
$data1[28,1] = '06:46'
$data1[28,2] = '215C_EL PATO'
$data1[28,3] = 110
$data1[28,4] = 'LP1912'
$data1[29,0] = '05:26:08'
$data1[29,1] = '06:47'
$data1[29,2] = '215C_EL PATO'
$data1[29,3] = 81
$data1[29,4] = 'LP1912'
$data1[30,0] = '05:55:25'
$data1[30,1] = '06:59'
$data1[30,2] = '14_ABASTO'
$data1[30,3] = 64
$data1[30,4] = 'LP1912'
$data1[31,0] = '05:26:08'
$data1[31,1] = '07:00'
$data1[31,2] = '14_ABASTO'
$data1[31,3] = 94
$data1[31,4] = 'LP1912'
$data1[32,0] = '06:25:43'
$data1[32,1] = '07:01'
$data1[32,2] = '16_SANTA ANA'
$data1[32,3] = 36
$data1[32,4] = 'LP1912'
$data1[33,0] = '05:55:25'
$data1[33,1] = '07:04'
$data1[33,2] = '23_HERNANDEZ'
$data1[33,3] = 69
$data1[33,4] = 'LP1912'
$data1[34,0] = '05:26:08'
$data1[34,1] = '07:05'
$data1[34,2] = '23_HERNANDEZ'
$data1[34,3] = 99
$data1[34,4] = 'LP1912'
$data1[35,0] = '05:26:08'
$data1[35,1] = '07:05'
$data1[35,2] = '15_ABASTO'
$data1[35,3] = 99
$data1[35,4] = 'LP1912'
$data1[36,0] = '05:26:08'
$data1[36,1] = '07:06'
$data1[36,2] = '10_OLMOS'
$data1[36,3] = 100
$data1[36,4] = 'LP1912'
$data1[37,0] = '05:26:08'
$data1[37,1] = '07:07'
$data1[37,2] = '225_GOMEZ'
$data1[37,3] = 101
$data1[37,4] = 'LP1912'
$data1[38,0] = '05:26:08'
$data1[38,1] = '07:11'
$data1[38,2] = '215A_EL PATO'
$data1[38,3] = 105
$data1[38,4] = 'LP1912'
$data1[39,0] = '06:55:02'
$data1[39,1] = '07:12'
$data1[39,2] = '215A_EL PATO'
$data1[39,3] = 17
$data1[39,4] = 'LP1912'
$data1[40,0] = '06:25:43'
$data1[40,1] = '07:14'
$data1[40,2] = '26_HERNANDEZ'
$data1[40,3] = 49
$data1[40,4] = 'LP1912'
$data1[41,0] = '05:55:25'
$data1[41,1] = '07:15'
$data1[41,2] = '11_ETCHEVERRY'
$data1[41,3] = 80
$data1[41,4] = 'LP1912'
$data1[42,0] = '05:26:08'
$data1[42,1] = '07:16'
$data1[42,2] = '11_ETCHEVERRY'
$data1[42,3] = 110
$data1[42,4] = 'LP1912'
$data1[43,0] = '06:55:02'
$data1[43,1] = '07:17'
$data1[43,2] = '16_SANTA ANA'
$data1[43,3] = 22
$data1[43,4] = 'LP1912'
$data1[44,0] = '05:26:08'
$data1[44,1] = '07:21'
$data1[44,2] = '26_HERNANDEZ'
$data1[44,3] = 115
$data1[44,4] = 'LP1912'
$data1[45,0] = '05:26:08'
$data1[45,1] = '07:23'
$data1[45,2] = '10_OLMOS'
$data1[45,3] = 117
$data1[45,4] = 'LP1912'
$data1[46,0] = '05:55:25'
$data1[46,1] = '07:30'
$data1[46,2] = '10_OLMOS'
$data1[46,3] = 95
$data1[46,4] = 'LP1912'
$data1[47,0] = '05:55:25'
$data1[47,1] = '07:31'
$data1[47,2] = '11_ETCHEVERRY'
$data1[47,3] = 96
$data1[47,4] = 'LP1912'
$data1[48,0] = '05:55:25'
$data1[48,1] = '07:31'
$data1[48,2] = '16_SANTA ANA'
$data1[48,3] = 96
$data1[48,4] = 'LP1912'
$data1[49,0] = '06:55:02'
$data1[49,1] = '07:32'
$data1[49,2] = '16_SANTA ANA'
$data1[49,3] = 37
$data1[49,4] = 'LP1912'
$data1[50,0] = '06:55:02'
$data1[50,1] = '07:32'
$data1[50,2] = '11_ETCHEVERRY'
$data1[50,3] = 37
$data1[50,4] = 'LP1912'
$data1[51,0] = '05:55:25'
$data1[51,1] = '07:32'
$data1[51,2] = '84_COLONIA URQUIZA-ESC 49'
$data1[51,3] = 97
$data1[51,4] = 'LP1912'
$data1[52,0] = '07:19:29'
$data1[52,1] = '07:35'
$data1[52,2] = '23_HERNANDEZ'
$data1[52,3] = 16
$data1[52,4] = 'LP1912'
$data1[53,0] = '05:55:25'
$data1[53,1] = '07:36'
$data1[53,2] = '27_EL RETIRO'
$data1[53,3] = 101
$data1[53,4] = 'LP1912'
$data1[54,0] = '06:55:02'
$data1[54,1] = '07:37'
$data1[54,2] = '27_EL RETIRO'
$data1[54,3] = 42
$data1[54,4] = 'LP1912'
$data1[55,0] = '05:55:25'
$data1[55,1] = '07:39'
$data1[55,2] = '10_OLMOS'
$data1[55,3] = 104
$data1[55,4] = 'LP1912'
$data1[56,0] = '07:19:29'
$data1[56,1] = '07:46'
$data1[56,2] = '16_SANTA ANA'
$data1[56,3] = 27
$data1[56,4] = 'LP1912'
$data1[57,0] = '05:55:25'
$data1[57,1] = '07:47'
$data1[57,2] = '14_ABASTO'
$data1[57,3] = 112
$data1[57,4] = 'LP1912'
$data1[58,0] = '06:55:02'
$data1[58,1] = '07:48'
$data1[58,2] = '14_ABASTO'
$data1[58,3] = 53
$data1[58,4] = 'LP1912'
$data1[59,0] = '07:50:16'
$data1[59,1] = '07:50'
$data1[59,2] = '10_OLMOS'
$data1[59,3] = 0
$data1[59,4] = 'LP1912'
$data1[60,0] = '05:55:25'
$data1[60,1] = '07:51'
$data1[60,2] = '215D_EL PATO'
$data1[60,3] = 116
$data1[60,4] = 'LP1912'
$data1[61,0] = '06:55:02'
$data1[61,1] = '07:52'
$data1[61,2] = '215D_EL PATO'
$data1[61,3] = 57
$data1[61,4] = 'LP1912'
$data1[62,0] = '07:19:29'
$data1[62,1] = '07:59'
$data1[62,2] = '23_HERNANDEZ'
$data1[62,3] = 40
$data1[62,4] = 'LP1912'
$data1[63,0] = '06:25:43'
$data1[63,1] = '08:01'
$data1[63,2] = '23_HERNANDEZ'
$data1[63,3] = 96
$data1[63,4] = 'LP1912'
$data1[64,0] = '06:55:02'
$data1[64,1] = '08:03'
$data1[64,2] = '23_HERNANDEZ'
$data1[64,3] = 68
$data1[64,4] = 'LP1912'
$data1[65,0] = '07:19:29'
$data1[65,1] = '08:03'
$data1[65,2] = '11_ETCHEVERRY'
$data1[65,3] = 44
$data1[65,4] = 'LP1912'
$data1[66,0] = '08:02:22'
$data1[66,1] = '08:05'
$data1[66,2] = '23_HERNANDEZ'
$data1[66,3] = 3
$data1[66,4] = 'LP1912'
$data1[67,0] = '07:19:29'
$data1[67,1] = '08:10'
$data1[67,2] = '16_SANTA ANA'
$data1[67,3] = 51
$data1[67,4] = 'LP1912'
$data1[68,0] = '07:50:16'
$data1[68,1] = '08:11'
$data1[68,2] = '16_SANTA ANA'
$data1[68,3] = 21
$data1[68,4] = 'LP1912'
$data1[69,0] = '06:25:43'
$data1[69,1] = '08:12'
$data1[69,2] = '15_ABASTO'
$data1[69,3] = 107
$data1[69,4] = 'LP1912'
$data1[70,0] = '07:50:16'
$data1[70,1] = '08:13'
$data1[70,2] = '10_OLMOS'
$data1[70,3] = 23
$data1[70,4] = 'LP1912'
$data1[71,0] = '06:55:02'
$data1[71,1] = '08:21'
$data1[71,2] = '26_HERNANDEZ'
$data1[71,3] = 86
$data1[71,4] = 'LP1912'
$data1[72,0] = '06:25:43'
$data1[72,1] = '08:22'
$data1[72,2] = '16_P MOR-SANTA ANA'
$data1[72,3] = 117
$data1[72,4] = 'LP1912'
$data1[73,0] = '06:55:02'
$data1[73,1] = '08:23'
$data1[73,2] = '16_P MOR-SANTA ANA'
$data1[73,3] = 88
$data1[73,4] = 'LP1912'
$data1[74,0] = '06:25:43'
$data1[74,1] = '08:23'
$data1[74,2] = '215B_EL PATO'
$data1[74,3] = 118
$data1[74,4] = 'LP1912'
$data1[75,0] = '06:55:02'
$data1[75,1] = '08:27'
$data1[75,2] = '84_COLONIA URQUIZA-ESC 49'
$data1[75,3] = 92
$data1[75,4] = 'LP1912'
$data1[76,0] = '07:50:16'
$data1[76,1] = '08:30'
$data1[76,2] = '23_HERNANDEZ'
$data1[76,3] = 40
$data1[76,4] = 'LP1912'
$data1[77,0] = '08:02:22'
$data1[77,1] = '08:33'
$data1[77,2] = '10_OLMOS'
$data1[77,3] = 31
$data1[77,4] = 'LP1912'
$data1[78,0] = '08:02:22'
$data1[78,1] = '08:34'
$data1[78,2] = '23_HERNANDEZ'
$data1[78,3] = 32
$data1[78,4] = 'LP1912'
$data1[79,0] = '08:32:09'
$data1[79,1] = '08:37'
$data1[79,2] = '23_HERNANDEZ'
$data1[79,3] = 5
$data1[79,4] = 'LP1912'
$data1[80,0] = '06:55:02'
$data1[80,1] = '08:42'
$data1[80,2] = '81_EL PELIGRO'
$data1[80,3] = 107
$data1[80,4] = 'LP1912'
$data1[81,0] = '07:19:29'
$data1[81,1] = '08:43'
$data1[81,2] = '14_ABASTO'
$data1[81,3] = 84
$data1[81,4] = 'LP1912'
$data1[82,0] = '07:50:16'
$data1[82,1] = '08:44'
$data1[82,2] = '14_ABASTO'
$data1[82,3] = 54
$data1[82,4] = 'LP1912'
$data1[83,0] = '08:32:09'
$data1[83,1] = '08:53'
$data1[83,2] = '10_OLMOS'
$data1[83,3] = 21
$data1[83,4] = 'LP1912'
$data1[84,0] = '06:55:02'
$data1[84,1] = '08:54'
$data1[84,2] = '17_ROMERO'
$data1[84,3] = 119
$data1[84,4] = 'LP1912'
$data1[85,0] = '07:19:29'
$data1[85,1] = '09:01'
$data1[85,2] = '215A_EL PATO'
$data1[85,3] = 102
$data1[85,4] = 'LP1912'
$data1[86,0] = '07:50:16'
$data1[86,1] = '09:02'
$data1[86,2] = '215A_EL PATO'
$data1[86,3] = 72
$data1[86,4] = 'LP1912'
$data1[87,0] = '08:02:22'
$data1[87,1] = '09:03'
$data1[87,2] = '11_ETCHEVERRY'
$data1[87,3] = 61
$data1[87,4] = 'LP1912'
$data1[88,0] = '08:32:09'
$data1[88,1] = '09:04'
$data1[88,2] = '11_ETCHEVERRY'
$data1[88,3] = 32
$data1[88,4] = 'LP1912'
$data1[89,0] = '08:32:09'
$data1[89,1] = '09:05'
$data1[89,2] = '23_HERNANDEZ'
$data1[89,3] = 33
$data1[89,4] = 'LP1912'
$data1[90,0] = '07:19:29'
$data1[90,1] = '09:10'
$data1[90,2] = '16_P MOR-SANTA ANA'
$data1[90,3] = 111
$data1[90,4] = 'LP1912'
$data1[91,0] = '07:50:16'
$data1[91,1] = '09:11'
$data1[91,2] = '16_P MOR-SANTA ANA'
$data1[91,3] = 81
$data1[91,4] = 'LP1912'
$data1[92,0] = '08:32:09'
$data1[92,1] = '09:13'
$data1[92,2] = '10_OLMOS'
$data1[92,3] = 41
$data1[92,4] = 'LP1912'
$data1[93,0] = '07:19:29'
$data1[93,1] = '09:16'
$data1[93,2] = '27_EL RETIRO'
$data1[93,3] = 117
$data1[93,4] = 'LP1912'
$data1[94,0] = '07:50:16'
$data1[94,1] = '09:17'
$data1[94,2] = '27_EL RETIRO'
$data1[94,3] = 87
$data1[94,4] = 'LP1912'
$data1[95,0] = '07:50:16'
$data1[95,1] = '09:21'
$data1[95,2] = '26_HERNANDEZ'
$data1[95,3] = 91
$data1[95,4] = 'LP1912'
$data1[96,0] = '08:02:22'
$data1[96,1] = '09:22'
$data1[96,2] = '16_SANTA ANA'
$data1[96,3] = 80
$data1[96,4] = 'LP1912'
$data1[97,0] = '08:02:22'
$data1[97,1] = '09:23'
$data1[97,2] = '11_ETCHEVERRY'
$data1[97,3] = 81
$data1[97,4] = 'LP1912'
$data1[98,0] = '07:50:16'
$data1[98,1] = '09:23'
$data1[98,2] = '17_ROMERO'
$data1[98,3] = 93
$data1[98,4] = 'LP1912'
$data1[99,0] = '08:32:09'
$data1[99,1] = '09:23'
$data1[99,2] = '16_SANTA ANA'
$data1[99,3] = 51
$data1[99,4] = 'LP1912'
$data1[100,0] = '07:50:16'
$data1[100,1] = '09:24'
$data1[100,2] = '11_ETCHEVERRY'
$data1[100,3] = 94
$data1[100,4] = 'LP1912'
$data1[101,0] = '07:50:16'
$data1[101,1] = '09:28'
$data1[101,2] = '16_SANTA ANA'
$data1[101,3] = 98
$data1[101,4] = 'LP1912'
$data1[102,0] = '07:50:16'
$data1[102,1] = '09:32'
$data1[102,2] = '15_ABASTO'
$data1[102,3] = 102
$data1[102,4] = 'LP1912'
$data1[103,0] = '07:50:16'
$data1[103,1] = '09:33'
$data1[103,2] = '10_OLMOS'
$data1[103,3] = 103
$data1[103,4] = 'LP1912'
$data1[104,0] = '08:56:29'
$data1[104,1] = '09:34'
$data1[104,2] = '16_SANTA ANA'
$data1[104,3] = 38
$data1[104,4] = 'LP1912'
$data1[105,0] = '08:56:29'
$data1[105,1] = '09:34'
$data1[105,2] = '23_HERNANDEZ'
$data1[105,3] = 38
$data1[105,4] = 'LP1912'
$data1[106,0] = '08:48:08'
$data1[106,1] = '09:35'
$data1[106,2] = '23_HERNANDEZ'
$data1[106,3] = 47
$data1[106,4] = 'LP1912'
$data1[107,0] = '08:32:09'
$data1[107,1] = '09:35'
$data1[107,2] = '16_SANTA ANA'
$data1[107,3] = 63
$data1[107,4] = 'LP1912'
$data1[108,0] = '09:35:26'
$data1[108,1] = '09:39'
$data1[108,2] = '23_HERNANDEZ'
$data1[108,3] = 4
$data1[108,4] = 'LP1912'
$data1[109,0] = '07:50:16'
$data1[109,1] = '09:42'
$data1[109,2] = '215C_EL PATO'
$data1[109,3] = 112
$data1[109,4] = 'LP1912'
$data1[110,0] = '08:02:22'
$data1[110,1] = '09:43'
$data1[110,2] = '14_ABASTO'
$data1[110,3] = 101
$data1[110,4] = 'LP1912'
$data1[111,0] = '07:50:16'
$data1[111,1] = '09:44'
$data1[111,2] = '14_ABASTO'
$data1[111,3] = 114
$data1[111,4] = 'LP1912'
$data1[112,0] = '09:35:26'
$data1[112,1] = '09:46'
$data1[112,2] = '16_SANTA ANA'
$data1[112,3] = 11
$data1[112,4] = 'LP1912'
$data1[113,0] = '08:32:09'
$data1[113,1] = '09:52'
$data1[113,2] = '15_ABASTO'
$data1[113,3] = 80
$data1[113,4] = 'LP1912'
$data1[114,0] = '08:56:29'
$data1[114,1] = '09:53'
$data1[114,2] = '10_OLMOS'
$data1[114,3] = 57
$data1[114,4] = 'LP1912'
$data1[115,0] = '09:35:26'
$data1[115,1] = '09:58'
$data1[115,2] = '16_SANTA ANA'
$data1[115,3] = 23
$data1[115,4] = 'LP1912'
$data1[116,0] = '09:35:26'
$data1[116,1] = '10:03'
$data1[116,2] = '11_ETCHEVERRY'
$data1[116,3] = 28
$data1[116,4] = 'LP1912'
$data1[117,0] = '08:56:29'
$data1[117,1] = '10:10'
$data1[117,2] = '16_P MOR-SANTA ANA'
$data1[117,3] = 74
$data1[117,4] = 'LP1912'
$data1[118,0] = '08:32:09'
$data1[118,1] = '10:11'
$data1[118,2] = '16_P MOR-SANTA ANA'
$data1[118,3] = 99
$data1[118,4] = 'LP1912'
$data1[119,0] = '09:35:26'
$data1[119,1] = '10:12'
$data1[119,2] = '15_ABASTO'
$data1[119,3] = 37
$data1[119,4] = 'LP1912'
$data1[120,0] = '09:35:26'
$data1[120,1] = '10:13'
$data1[120,2] = '10_OLMOS'
$data1[120,3] = 38
$data1[120,4] = 'LP1912'
$data1[121,0] = '08:32:09'
$data1[121,1] = '10:21'
$data1[121,2] = '26_HERNANDEZ'
$data1[121,3] = 109
$data1[121,4] = 'LP1912'
$data1[122,0] = '08:32:09'
$data1[122,1] = '10:22'
$data1[122,2] = '17_ROMERO'
$data1[122,3] = 110
$data1[122,4] = 'LP1912'
$data1[123,0] = '09:35:26'
$data1[123,1] = '10:23'
$data1[123,2] = '11_ETCHEVERRY'
$data1[123,3] = 48
$data1[123,4] = 'LP1912'
$data1[124,0] = '08:56:29'
$data1[124,1] = '10:26'
$data1[124,2] = '215A_EL PATO'
$data1[124,3] = 90
$data1[124,4] = 'LP1912'
$data1[125,0] = '08:32:09'
$data1[125,1] = '10:27'
$data1[125,2] = '215A_EL PATO'
$data1[125,3] = 115
$data1[125,4] = 'LP1912'
$data1[126,0] = '10:29:57'
$data1[126,1] = '10:29'
$data1[126,2] = '16_SANTA ANA'
$data1[126,3] = 0
$data1[126,4] = 'LP1912'
$data1[127,0] = '10:29:57'
$data1[127,1] = '10:31'
$data1[127,2] = '10_OLMOS'
$data1[127,3] = 2
$data1[127,4] = 'LP1912'
$data1[128,0] = '10:29:57'
$data1[128,1] = '10:34'
$data1[128,2] = '16_SANTA ANA'
$data1[128,3] = 5
$data1[128,4] = 'LP1912'
$data1[129,0] = '09:35:26'
$data1[129,1] = '10:34'
$data1[129,2] = '23_HERNANDEZ'
$data1[129,3] = 59
$data1[129,4] = 'LP1912'
$data1[130,0] = '10:29:57'
$data1[130,1] = '10:39'
$data1[130,2] = '23_HERNANDEZ'
$data1[130,3] = 10
$data1[130,4] = 'LP1912'
$data1[131,0] = '10:29:57'
$data1[131,1] = '10:41'
$data1[131,2] = '17_ROMERO'
$data1[131,3] = 12
$data1[131,4] = 'LP1912'
$data1[132,0] = '08:48:08'
$data1[132,1] = '10:42'
$data1[132,2] = '17_ROMERO'
$data1[132,3] = 114
$data1[132,4] = 'LP1912'
$data1[133,0] = '08:56:29'
$data1[133,1] = '10:43'
$data1[133,2] = '14_ABASTO'
$data1[133,3] = 107
$data1[133,4] = 'LP1912'
$data1[134,0] = '08:48:08'
$data1[134,1] = '10:44'
$data1[134,2] = '14_ABASTO'
$data1[134,3] = 116
$data1[134,4] = 'LP1912'
$data1[135,0] = '10:29:57'
$data1[135,1] = '10:51'
$data1[135,2] = '15_ABASTO'
$data1[135,3] = 22
$data1[135,4] = 'LP1912'
$data1[136,0] = '10:29:57'
$data1[136,1] = '10:52'
$data1[136,2] = '10_OLMOS'
$data1[136,3] = 23
$data1[136,4] = 'LP1912'
$data1[137,0] = '09:35:26'
$data1[137,1] = '10:54'
$data1[137,2] = '27_EL RETIRO'
$data1[137,3] = 79
$data1[137,4] = 'LP1912'
$data1[138,0] = '10:29:57'
$data1[138,1] = '10:56'
$data1[138,2] = '27_EL RETIRO'
$data1[138,3] = 27
$data1[138,4] = 'LP1912'
$data1[139,0] = '10:59:49'
$data1[139,1] = '10:59'
$data1[139,2] = '16_SANTA ANA'
$data1[139,3] = 0
$data1[139,4] = 'LP1912'
$data1[140,0] = '10:29:57'
$data1[140,1] = '11:01'
$data1[140,2] = '215C_EL PATO'
$data1[140,3] = 32
$data1[140,4] = 'LP1912'
$data1[141,0] = '09:35:26'
$data1[141,1] = '11:02'
$data1[141,2] = '215C_EL PATO'
$data1[141,3] = 87
$data1[141,4] = 'LP1912'
$data1[142,0] = '10:29:57'
$data1[142,1] = '11:03'
$data1[142,2] = '11_ETCHEVERRY'
$data1[142,3] = 34
$data1[142,4] = 'LP1912'
$data1[143,0] = '10:29:57'
$data1[143,1] = '11:04'
$data1[143,2] = '23_HERNANDEZ'
$data1[143,3] = 35
$data1[143,4] = 'LP1912'
$data1[144,0] = '09:35:26'
$data1[144,1] = '11:06'
$data1[144,2] = '16_P MOR-167 Y 521'
$data1[144,3] = 91
$data1[144,4] = 'LP1912'
$data1[145,0] = '10:59:49'
$data1[145,1] = '11:06'
$data1[145,2] = '23_HERNANDEZ'
$data1[145,3] = 7
$data1[145,4] = 'LP1912'
$data1[146,0] = '10:29:57'
$data1[146,1] = '11:11'
$data1[146,2] = '15_ABASTO'
$data1[146,3] = 42
$data1[146,4] = 'LP1912'
$data1[147,0] = '10:59:49'
$data1[147,1] = '11:11'
$data1[147,2] = '10_OLMOS'
$data1[147,3] = 12
$data1[147,4] = 'LP1912'
$data1[148,0] = '10:59:49'
$data1[148,1] = '11:12'
$data1[148,2] = '15_ABASTO'
$data1[148,3] = 13
$data1[148,4] = 'LP1912'
$data1[149,0] = '09:35:26'
$data1[149,1] = '11:19'
$data1[149,2] = '86_EST CHICA-ESC AGRARIA'
$data1[149,3] = 104
$data1[149,4] = 'LP1912'
$data1[150,0] = '10:29:57'
$data1[150,1] = '11:20'
$data1[150,2] = '26_HERNANDEZ'
$data1[150,3] = 51
$data1[150,4] = 'LP1912'
$data1[151,0] = '09:35:26'
$data1[151,1] = '11:21'
$data1[151,2] = '26_HERNANDEZ'
$data1[151,3] = 106
$data1[151,4] = 'LP1912'
$data1[152,0] = '10:29:57'
$data1[152,1] = '11:26'
$data1[152,2] = '225_C ROCA-H SUR'
$data1[152,3] = 57
$data1[152,4] = 'LP1912'
$data1[153,0] = '09:35:26'
$data1[153,1] = '11:27'
$data1[153,2] = '225_C ROCA-H SUR'
$data1[153,3] = 112
$data1[153,4] = 'LP1912'
$data1[154,0] = '11:30:45'
$data1[154,1] = '11:30'
$data1[154,2] = '16_SANTA ANA'
$data1[154,3] = 0
$data1[154,4] = 'LP1912'
$data1[155,0] = '10:29:57'
$data1[155,1] = '11:31'
$data1[155,2] = '81_EL PELIGRO'
$data1[155,3] = 62
$data1[155,4] = 'LP1912'
$data1[156,0] = '11:30:45'
$data1[156,1] = '11:31'
$data1[156,2] = '16_SANTA ANA'
$data1[156,3] = 1
$data1[156,4] = 'LP1912'
$data1[157,0] = '09:35:26'
$data1[157,1] = '11:32'
$data1[157,2] = '81_EL PELIGRO'
$data1[157,3] = 117
$data1[157,4] = 'LP1912'
$data1[158,0] = '10:59:49'
$data1[158,1] = '11:34'
$data1[158,2] = '23_HERNANDEZ'
$data1[158,3] = 35
$data1[158,4] = 'LP1912'
$data1[159,0] = '10:29:57'
$data1[159,1] = '11:35'
$data1[159,2] = '11_ETCHEVERRY'
$data1[159,3] = 66
$data1[159,4] = 'LP1912'
$data1[160,0] = '10:29:57'
$data1[160,1] = '11:40'
$data1[160,2] = '10_OLMOS'
$data1[160,3] = 71
$data1[160,4] = 'LP1912'
$data1[161,0] = '10:29:57'
$data1[161,1] = '11:41'
$data1[161,2] = '17_ROMERO'
$data1[161,3] = 72
$data1[161,4] = 'LP1912'
$data1[162,0] = '10:59:49'
$data1[162,1] = '11:42'
$data1[162,2] = '11_ETCHEVERRY'
$data1[162,3] = 43
$data1[162,4] = 'LP1912'
$data1[163,0] = '10:59:49'
$data1[163,1] = '11:43'
$data1[163,2] = '10_OLMOS'
$data1[163,3] = 44
$data1[163,4] = 'LP1912'
$data1[164,0] = '11:30:45'
$data1[164,1] = '11:44'
$data1[164,2] = '11_ETCHEVERRY'
$data1[164,3] = 14
$data1[164,4] = 'LP1912'
$data1[165,0] = '10:29:57'
$data1[165,1] = '11:50'
$data1[165,2] = '215B_EL PATO'
$data1[165,3] = 81
$data1[165,4] = 'LP1912'
$data1[166,0] = '10:59:49'
$data1[166,1] = '11:51'
$data1[166,2] = '215B_EL PATO'
$data1[166,3] = 52
$data1[166,4] = 'LP1912'
$data1[167,0] = '10:59:49'
$data1[167,1] = '11:52'
$data1[167,2] = '15_ABASTO'
$data1[167,3] = 53
$data1[167,4] = 'LP1912'
$data1[168,0] = '11:56:55'
$data1[168,1] = '11:56'
$data1[168,2] = '16_SANTA ANA'
$data1[168,3] = 0
$data1[168,4] = 'LP1912'
$data1[169,0] = '10:29:57'
$data1[169,1] = '11:58'
$data1[169,2] = '225_GOMEZ'
$data1[169,3] = 89
$data1[169,4] = 'LP1912'
$data1[170,0] = '11:30:45'
$data1[170,1] = '11:59'
$data1[170,2] = '225_GOMEZ'
$data1[170,3] = 29
$data1[170,4] = 'LP1912'
$data1[171,0] = '10:29:57'
$data1[171,1] = '12:01'
$data1[171,2] = '84_COLONIA URQUIZA-ESC 49'
$data1[171,3] = 92
$data1[171,4] = 'LP1912'
$data1[172,0] = '10:59:49'
$data1[172,1] = '12:02'
$data1[172,2] = '84_COLONIA URQUIZA-ESC 49'
$data1[172,3] = 63
$data1[172,4] = 'LP1912'
$data1[173,0] = '11:30:45'
$data1[173,1] = '12:04'
$data1[173,2] = '23_HERNANDEZ'
$data1[173,3] = 34
$data1[173,4] = 'LP1912'
$data1[174,0] = '10:29:57'
$data1[174,1] = '12:06'
$data1[174,2] = '16_P MOR-SANTA ANA'
$data1[174,3] = 97
$data1[174,4] = 'LP1912'
$data1[175,0] = '10:59:49'
$data1[175,1] = '12:06'
$data1[175,2] = '14_ABASTO'
$data1[175,3] = 67
$data1[175,4] = 'LP1912'
$data1[176,0] = '10:59:49'
$data1[176,1] = '12:10'
$data1[176,2] = '10_OLMOS'
$data1[176,3] = 71
$data1[176,4] = 'LP1912'
$data1[177,0] = '11:56:55'
$data1[177,1] = '12:12'
$data1[177,2] = '10_OLMOS'
$data1[177,3] = 16
$data1[177,4] = 'LP1912'
$data1[178,0] = '10:29:57'
$data1[178,1] = '12:13'
$data1[178,2] = '17_ROMERO'
$data1[178,3] = 104
$data1[178,4] = 'LP1912'
$data1[179,0] = '11:56:55'
$data1[179,1] = '12:14'
$data1[179,2] = '17_ROMERO'
$data1[179,3] = 18
$data1[179,4] = 'LP1912'
$data1[180,0] = '10:29:57'
$data1[180,1] = '12:15'
$data1[180,2] = '14_ABASTO'
$data1[180,3] = 106
$data1[180,4] = 'LP1912'
$data1[181,0] = '10:59:49'
$data1[181,1] = '12:20'
$data1[181,2] = '14_ABASTO'
$data1[181,3] = 81
$data1[181,4] = 'LP1912'
$data1[182,0] = '10:29:57'
$data1[182,1] = '12:20'
$data1[182,2] = '26_HERNANDEZ'
$data1[182,3] = 111
$data1[182,4] = 'LP1912'
$data1[183,0] = '10:29:57'
$data1[183,1] = '12:20'
$data1[183,2] = '215A_EL PATO'
$data1[183,3] = 111
$data1[183,4] = 'LP1912'
$data1[184,0] = '10:59:49'
$data1[184,1] = '12:21'
$data1[184,2] = '26_HERNANDEZ'
$data1[184,3] = 82
$data1[184,4] = 'LP1912'
$data1[185,0] = '12:21:08'
$data1[185,1] = '12:21'
$data1[185,2] = '16_SANTA ANA'
$data1[185,3] = 0
$data1[185,4] = 'LP1912'
$data1[186,0] = '12:21:08'
$data1[186,1] = '12:21'
$data1[186,2] = '215A_EL PATO'
$data1[186,3] = 0
$data1[186,4] = 'LP1912'
$data1[187,0] = '10:59:49'
$data1[187,1] = '12:30'
$data1[187,2] = '17_ROMERO'
$data1[187,3] = 91
$data1[187,4] = 'LP1912'
$data1[188,0] = '11:56:55'
$data1[188,1] = '12:34'
$data1[188,2] = '11_ETCHEVERRY'
$data1[188,3] = 38
$data1[188,4] = 'LP1912'
$data1[189,0] = '11:56:55'
$data1[189,1] = '12:34'
$data1[189,2] = '23_HERNANDEZ'
$data1[189,3] = 38
$data1[189,4] = 'LP1912'
$data1[190,0] = '12:21:08'
$data1[190,1] = '12:35'
$data1[190,2] = '11_ETCHEVERRY'
$data1[190,3] = 14
$data1[190,4] = 'LP1912'
$data1[191,0] = '12:21:08'
$data1[191,1] = '12:35'
$data1[191,2] = '23_HERNANDEZ'
$data1[191,3] = 14
$data1[191,4] = 'LP1912'
$data1[192,0] = '10:59:49'
$data1[192,1] = '12:36'
$data1[192,2] = '27_EL RETIRO'
$data1[192,3] = 97
$data1[192,4] = 'LP1912'
$data1[193,0] = '12:21:08'
$data1[193,1] = '12:37'
$data1[193,2] = '27_EL RETIRO'
$data1[193,3] = 16
$data1[193,4] = 'LP1912'
$data1[194,0] = '10:59:49'
$data1[194,1] = '12:38'
$data1[194,2] = '17_179 Y 38'
$data1[194,3] = 99
$data1[194,4] = 'LP1912'
$data1[195,0] = '11:56:55'
$data1[195,1] = '12:40'
$data1[195,2] = '10_OLMOS'
$data1[195,3] = 44
$data1[195,4] = 'LP1912'
$data1[196,0] = '11:30:45'
$data1[196,1] = '12:41'
$data1[196,2] = '10_OLMOS'
$data1[196,3] = 71
$data1[196,4] = 'LP1912'
$data1[197,0] = '12:47:27'
$data1[197,1] = '12:47'
$data1[197,2] = '16_SANTA ANA'
$data1[197,3] = 0
$data1[197,4] = 'LP1912'
$data1[198,0] = '10:59:49'
$data1[198,1] = '12:48'
$data1[198,2] = '11_ETCHEVERRY'
$data1[198,3] = 109
$data1[198,4] = 'LP1912'
$data1[199,0] = '12:47:27'
$data1[199,1] = '12:48'
$data1[199,2] = '16_SANTA ANA'
$data1[199,3] = 1
$data1[199,4] = 'LP1912'
$data1[200,0] = '12:21:08'
$data1[200,1] = '12:49'
$data1[200,2] = '11_ETCHEVERRY'
$data1[200,3] = 28
$data1[200,4] = 'LP1912'
$data1[201,0] = '12:21:08'
$data1[201,1] = '12:55'
$data1[201,2] = '10_OLMOS'
$data1[201,3] = 34
$data1[201,4] = 'LP1912'
$data1[202,0] = '12:59:47'
$data1[202,1] = '13:00'
$data1[202,2] = '16_SANTA ANA'
$data1[202,3] = 1
$data1[202,4] = 'LP1912'
$data1[203,0] = '11:30:45'
$data1[203,1] = '13:01'
$data1[203,2] = '17_ROMERO'
$data1[203,3] = 91
$data1[203,4] = 'LP1912'
$data1[204,0] = '12:47:27'
$data1[204,1] = '13:02'
$data1[204,2] = '15_ABASTO'
$data1[204,3] = 15
$data1[204,4] = 'LP1912'
$data1[205,0] = '12:21:08'
$data1[205,1] = '13:03'
$data1[205,2] = '14_ABASTO'
$data1[205,3] = 42
$data1[205,4] = 'LP1912'
$data1[206,0] = '12:47:27'
$data1[206,1] = '13:04'
$data1[206,2] = '23_HERNANDEZ'
$data1[206,3] = 17
$data1[206,4] = 'LP1912'
$data1[207,0] = '12:59:47'
$data1[207,1] = '13:05'
$data1[207,2] = '23_HERNANDEZ'
$data1[207,3] = 6
$data1[207,4] = 'LP1912'
$data1[208,0] = '11:30:45'
$data1[208,1] = '13:06'
$data1[208,2] = '16_P MOR-SANTA ANA'
$data1[208,3] = 96
$data1[208,4] = 'LP1912'
$data1[209,0] = '12:21:08'
$data1[209,1] = '13:07'
$data1[209,2] = '16_P MOR-SANTA ANA'
$data1[209,3] = 46
$data1[209,4] = 'LP1912'
$data1[210,0] = '11:30:45'
$data1[210,1] = '13:07'
$data1[210,2] = '10_OLMOS'
$data1[210,3] = 97
$data1[210,4] = 'LP1912'
$data1[211,0] = '12:21:08'
$data1[211,1] = '13:08'
$data1[211,2] = '10_OLMOS'
$data1[211,3] = 47
$data1[211,4] = 'LP1912'
$data1[212,0] = '11:30:45'
$data1[212,1] = '13:13'
$data1[212,2] = '215D_EL PATO'
$data1[212,3] = 103
$data1[212,4] = 'LP1912'
$data1[213,0] = '12:21:08'
$data1[213,1] = '13:14'
$data1[213,2] = '215D_EL PATO'
$data1[213,3] = 53
$data1[213,4] = 'LP1912'
$data1[214,0] = '12:47:27'
$data1[214,1] = '13:14'
$data1[214,2] = '11_ETCHEVERRY'
$data1[214,3] = 27
$data1[214,4] = 'LP1912'
$data1[215,0] = '11:56:55'
$data1[215,1] = '13:20'
$data1[215,2] = '26_HERNANDEZ'
$data1[215,3] = 84
$data1[215,4] = 'LP1912'
$data1[216,0] = '11:30:45'
$data1[216,1] = '13:21'
$data1[216,2] = '26_HERNANDEZ'
$data1[216,3] = 111
$data1[216,4] = 'LP1912'
$data1[217,0] = '11:30:45'
$data1[217,1] = '13:25'
$data1[217,2] = '10_OLMOS'
$data1[217,3] = 115
$data1[217,4] = 'LP1912'
$data1[218,0] = '11:30:45'
$data1[218,1] = '13:26'
$data1[218,2] = '14_ABASTO'
$data1[218,3] = 116
$data1[218,4] = 'LP1912'
$data1[219,0] = '11:30:45'
$data1[219,1] = '13:26'
$data1[219,2] = '15_ABASTO'
$data1[219,3] = 116
$data1[219,4] = 'LP1912'
$data1[220,0] = '11:56:55'
$data1[220,1] = '13:27'
$data1[220,2] = '10_OLMOS'
$data1[220,3] = 91
$data1[220,4] = 'LP1912'
$data1[221,0] = '12:21:08'
$data1[221,1] = '13:27'
$data1[221,2] = '14_ABASTO'
$data1[221,3] = 66
$data1[221,4] = 'LP1912'
$data1[222,0] = '12:21:08'
$data1[222,1] = '13:28'
$data1[222,2] = '10_OLMOS'
$data1[222,3] = 67
$data1[222,4] = 'LP1912'
$data1[223,0] = '12:47:27'
$data1[223,1] = '13:31'
$data1[223,2] = '10_OLMOS'
$data1[223,3] = 44
$data1[223,4] = 'LP1912'
$data1[224,0] = '12:47:27'
$data1[224,1] = '13:32'
$data1[224,2] = '10_OLMOS'
$data1[224,3] = 45
$data1[224,4] = 'LP1912'
$data1[225,0] = '12:59:47'
$data1[225,1] = '13:33'
$data1[225,2] = '10_OLMOS'
$data1[225,3] = 34
$data1[225,4] = 'LP1912'
$data1[226,0] = '13:33:42'
$data1[226,1] = '13:33'
$data1[226,2] = '16_SANTA ANA'
$data1[226,3] = 0
$data1[226,4] = 'LP1912'
$data1[227,0] = '13:33:42'
$data1[227,1] = '13:34'
$data1[227,2] = '16_SANTA ANA'
$data1[227,3] = 1
$data1[227,4] = 'LP1912'
$data1[228,0] = '13:33:42'
$data1[228,1] = '13:34'
$data1[228,2] = '23_HERNANDEZ'
$data1[228,3] = 1
$data1[228,4] = 'LP1912'
$data1[229,0] = '11:56:55'
$data1[229,1] = '13:36'
$data1[229,2] = '15_ABASTO'
$data1[229,3] = 100
$data1[229,4] = 'LP1912'
$data1[230,0] = '13:33:42'
$data1[230,1] = '13:38'
$data1[230,2] = '14_ABASTO'
$data1[230,3] = 5
$data1[230,4] = 'LP1912'
$data1[231,0] = '11:56:55'
$data1[231,1] = '13:46'
$data1[231,2] = '17_ROMERO'
$data1[231,3] = 110
$data1[231,4] = 'LP1912'
$data1[232,0] = '12:59:47'
$data1[232,1] = '13:50'
$data1[232,2] = '11_ETCHEVERRY'
$data1[232,3] = 51
$data1[232,4] = 'LP1912'
$data1[233,0] = '11:56:55'
$data1[233,1] = '13:50'
$data1[233,2] = '215A_EL PATO'
$data1[233,3] = 114
$data1[233,4] = 'LP1912'
$data1[234,0] = '12:21:08'
$data1[234,1] = '13:51'
$data1[234,2] = '215A_EL PATO'
$data1[234,3] = 90
$data1[234,4] = 'LP1912'
$data1[235,0] = '11:56:55'
$data1[235,1] = '13:55'
$data1[235,2] = '225_GOMEZ'
$data1[235,3] = 119
$data1[235,4] = 'LP1912'
$data1[236,0] = '12:59:47'
$data1[236,1] = '13:56'
$data1[236,2] = '16_P MOR-167 Y 521'
$data1[236,3] = 57
$data1[236,4] = 'LP1912'
$data1[237,0] = '12:21:08'
$data1[237,1] = '13:56'
$data1[237,2] = '225_GOMEZ'
$data1[237,3] = 95
$data1[237,4] = 'LP1912'
$data1[238,0] = '12:47:27'
$data1[238,1] = '13:58'
$data1[238,2] = '16_P MOR-167 Y 521'
$data1[238,3] = 71
$data1[238,4] = 'LP1912'
$data1[239,0] = '13:59:06'
$data1[239,1] = '13:59'
$data1[239,2] = '16_SANTA ANA'
$data1[239,3] = 0
$data1[239,4] = 'LP1912'
$data1[240,0] = '12:21:08'
$data1[240,1] = '14:00'
$data1[240,2] = '16_P MOR-167 Y 521'
$data1[240,3] = 99
$data1[240,4] = 'LP1912'
$data1[241,0] = '13:59:06'
$data1[241,1] = '14:00'
$data1[241,2] = '16_SANTA ANA'
$data1[241,3] = 1
$data1[241,4] = 'LP1912'
$data1[242,0] = '13:33:42'
$data1[242,1] = '14:04'
$data1[242,2] = '23_HERNANDEZ'
$data1[242,3] = 31
$data1[242,4] = 'LP1912'
$data1[243,0] = '12:21:08'
$data1[243,1] = '14:04'
$data1[243,2] = '17_ROMERO'
$data1[243,3] = 103
$data1[243,4] = 'LP1912'
$data1[244,0] = '13:59:06'
$data1[244,1] = '14:05'
$data1[244,2] = '23_HERNANDEZ'
$data1[244,3] = 6
$data1[244,4] = 'LP1912'
$data1[245,0] = '12:21:08'
$data1[245,1] = '14:08'
$data1[245,2] = '23_HERNANDEZ'
$data1[245,3] = 107
$data1[245,4] = 'LP1912'
$data1[246,0] = '12:59:47'
$data1[246,1] = '14:11'
$data1[246,2] = '23_HERNANDEZ'
$data1[246,3] = 72
$data1[246,4] = 'LP1912'
$data1[247,0] = '13:33:42'
$data1[247,1] = '14:12'
$data1[247,2] = '15_ABASTO'
$data1[247,3] = 39
$data1[247,4] = 'LP1912'
$data1[248,0] = '12:47:27'
$data1[248,1] = '14:16'
$data1[248,2] = '27_EL RETIRO'
$data1[248,3] = 89
$data1[248,4] = 'LP1912'
$data1[249,0] = '12:21:08'
$data1[249,1] = '14:17'
$data1[249,2] = '27_EL RETIRO'
$data1[249,3] = 116
$data1[249,4] = 'LP1912'
$data1[250,0] = '12:59:47'
$data1[250,1] = '14:19'
$data1[250,2] = '215C_EL PATO'
$data1[250,3] = 80
$data1[250,4] = 'LP1912'
$data1[251,0] = '12:21:08'
$data1[251,1] = '14:20'
$data1[251,2] = '215C_EL PATO'
$data1[251,3] = 119
$data1[251,4] = 'LP1912'
$data1[252,0] = '12:47:27'
$data1[252,1] = '14:21'
$data1[252,2] = '26_HERNANDEZ'
$data1[252,3] = 94
$data1[252,4] = 'LP1912'
$data1[253,0] = '14:24:16'
$data1[253,1] = '14:25'
$data1[253,2] = '16_SANTA ANA'
$data1[253,3] = 1
$data1[253,4] = 'LP1912'
$data1[254,0] = '13:59:06'
$data1[254,1] = '14:28'
$data1[254,2] = '15_ABASTO'
$data1[254,3] = 29
$data1[254,4] = 'LP1912'
$data1[255,0] = '14:24:16'
$data1[255,1] = '14:35'
$data1[255,2] = '23_HERNANDEZ'
$data1[255,3] = 11
$data1[255,4] = 'LP1912'
$data1[256,0] = '14:24:16'
$data1[256,1] = '14:44'
$data1[256,2] = '15_ABASTO'
$data1[256,3] = 20
$data1[256,4] = 'LP1912'
$data1[257,0] = '13:33:42'
$data1[257,1] = '14:44'
$data1[257,2] = '14_ABASTO'
$data1[257,3] = 71
$data1[257,4] = 'LP1912'
$data1[258,0] = '14:45:17'
$data1[258,1] = '14:45'
$data1[258,2] = '15_ABASTO'
$data1[258,3] = 0
$data1[258,4] = 'LP1912'
$data1[259,0] = '12:47:27'
$data1[259,1] = '14:45'
$data1[259,2] = '14_ABASTO'
$data1[259,3] = 118
$data1[259,4] = 'LP1912'
$data1[260,0] = '14:45:17'
$data1[260,1] = '14:46'
$data1[260,2] = '16_SANTA ANA'
$data1[260,3] = 1
$data1[260,4] = 'LP1912'
$data1[261,0] = '12:59:47'
$data1[261,1] = '14:56'
$data1[261,2] = '16_P MOR-SANTA ANA'
$data1[261,3] = 117
$data1[261,4] = 'LP1912'
$data1[262,0] = '14:56:20'
$data1[262,1] = '14:56'
$data1[262,2] = '16_SANTA ANA'
$data1[262,3] = 0
$data1[262,4] = 'LP1912'
$data1[263,0] = '13:59:06'
$data1[263,1] = '14:57'
$data1[263,2] = '16_P MOR-SANTA ANA'
$data1[263,3] = 58
$data1[263,4] = 'LP1912'
$data1[264,0] = '12:59:47'
$data1[264,1] = '14:58'
$data1[264,2] = '215B_EL PATO'
$data1[264,3] = 119
$data1[264,4] = 'LP1912'
$data1[265,0] = '13:33:42'
$data1[265,1] = '15:00'
$data1[265,2] = '81_EL PELIGRO'
$data1[265,3] = 87
$data1[265,4] = 'LP1912'
$data1[266,0] = '14:45:17'
$data1[266,1] = '15:05'
$data1[266,2] = '23_HERNANDEZ'
$data1[266,3] = 20
$data1[266,4] = 'LP1912'
$data1[267,0] = '13:33:42'
$data1[267,1] = '15:05'
$data1[267,2] = '10_OLMOS'
$data1[267,3] = 92
$data1[267,4] = 'LP1912'
$data1[268,0] = '13:59:06'
$data1[268,1] = '15:10'
$data1[268,2] = '17_ROMERO'
$data1[268,3] = 71
$data1[268,4] = 'LP1912'
$data1[269,0] = '13:33:42'
$data1[269,1] = '15:13'
$data1[269,2] = '11_ETCHEVERRY'
$data1[269,3] = 100
$data1[269,4] = 'LP1912'
$data1[270,0] = '13:59:06'
$data1[270,1] = '15:14'
$data1[270,2] = '11_ETCHEVERRY'
$data1[270,3] = 75
$data1[270,4] = 'LP1912'
$data1[271,0] = '14:56:20'
$data1[271,1] = '15:17'
$data1[271,2] = '16_SANTA ANA'
$data1[271,3] = 21
$data1[271,4] = 'LP1912'
$data1[272,0] = '13:33:42'
$data1[272,1] = '15:17'
$data1[272,2] = '26_HERNANDEZ'
$data1[272,3] = 104
$data1[272,4] = 'LP1912'
$data1[273,0] = '13:59:06'
$data1[273,1] = '15:18'
$data1[273,2] = '26_HERNANDEZ'
$data1[273,3] = 79
$data1[273,4] = 'LP1912'
$data1[274,0] = '14:56:20'
$data1[274,1] = '15:20'
$data1[274,2] = '15_ABASTO'
$data1[274,3] = 24
$data1[274,4] = 'LP1912'
$data1[275,0] = '14:24:16'
$data1[275,1] = '15:21'
$data1[275,2] = '26_HERNANDEZ'
$data1[275,3] = 57
$data1[275,4] = 'LP1912'
$data1[276,0] = '15:22:17'
$data1[276,1] = '15:22'
$data1[276,2] = '16_SANTA ANA'
$data1[276,3] = 0
$data1[276,4] = 'LP1912'
$data1[277,0] = '15:22:17'
$data1[277,1] = '15:22'
$data1[277,2] = '26_HERNANDEZ'
$data1[277,3] = 0
$data1[277,4] = 'LP1912'
$data1[278,0] = '14:24:16'
$data1[278,1] = '15:32'
$data1[278,2] = '84_COLONIA URQUIZA-ESC 49'
$data1[278,3] = 68
$data1[278,4] = 'LP1912'
$data1[279,0] = '13:59:06'
$data1[279,1] = '15:35'
$data1[279,2] = '23_HERNANDEZ'
$data1[279,3] = 96
$data1[279,4] = 'LP1912'
$data1[280,0] = '13:59:06'
$data1[280,1] = '15:37'
$data1[280,2] = '10_OLMOS'
$data1[280,3] = 98
$data1[280,4] = 'LP1912'
$data1[281,0] = '14:56:20'
$data1[281,1] = '15:38'
$data1[281,2] = '10_OLMOS'
$data1[281,3] = 42
$data1[281,4] = 'LP1912'
$data1[282,0] = '14:24:16'
$data1[282,1] = '15:38'
$data1[282,2] = '23_HERNANDEZ'
$data1[282,3] = 74
$data1[282,4] = 'LP1912'
$data1[283,0] = '14:45:17'
$data1[283,1] = '15:38'
$data1[283,2] = '215A_EL PATO'
$data1[283,3] = 53
$data1[283,4] = 'LP1912'
$data1[284,0] = '13:59:06'
$data1[284,1] = '15:39'
$data1[284,2] = '215A_EL PATO'
$data1[284,3] = 100
$data1[284,4] = 'LP1912'
$data1[285,0] = '14:56:20'
$data1[285,1] = '15:45'
$data1[285,2] = '14_ABASTO'
$data1[285,3] = 49
$data1[285,4] = 'LP1912'
$data1[286,0] = '14:56:20'
$data1[286,1] = '15:46'
$data1[286,2] = '16_P MOR-167 Y 521'
$data1[286,3] = 50
$data1[286,4] = 'LP1912'
$data1[287,0] = '14:24:16'
$data1[287,1] = '15:46'
$data1[287,2] = '14_ABASTO'
$data1[287,3] = 82
$data1[287,4] = 'LP1912'
$data1[288,0] = '13:59:06'
$data1[288,1] = '15:47'
$data1[288,2] = '16_P MOR-167 Y 521'
$data1[288,3] = 108
$data1[288,4] = 'LP1912'
$data1[289,0] = '13:59:06'
$data1[289,1] = '15:48'
$data1[289,2] = '14_ABASTO'
$data1[289,3] = 109
$data1[289,4] = 'LP1912'
$data1[290,0] = '14:56:20'
$data1[290,1] = '15:53'
$data1[290,2] = '11_ETCHEVERRY'
$data1[290,3] = 57
$data1[290,4] = 'LP1912'
$data1[291,0] = '13:59:06'
$data1[291,1] = '15:54'
$data1[291,2] = '11_ETCHEVERRY'
$data1[291,3] = 115
$data1[291,4] = 'LP1912'
$data1[292,0] = '15:53:28'
$data1[292,1] = '15:54'
$data1[292,2] = '11_ETCHEVERRY'
$data1[292,3] = 1
$data1[292,4] = 'LP1912'
$data1[293,0] = '15:22:17'
$data1[293,1] = '15:55'
$data1[293,2] = '16_SANTA ANA'
$data1[293,3] = 33
$data1[293,4] = 'LP1912'
$data1[294,0] = '15:53:28'
$data1[294,1] = '15:56'
$data1[294,2] = '11_ETCHEVERRY'
$data1[294,3] = 3
$data1[294,4] = 'LP1912'
$data1[295,0] = '14:24:16'
$data1[295,1] = '15:56'
$data1[295,2] = '17_ROMERO'
$data1[295,3] = 92
$data1[295,4] = 'LP1912'
$data1[296,0] = '13:59:06'
$data1[296,1] = '15:57'
$data1[296,2] = '27_EL RETIRO'
$data1[296,3] = 118
$data1[296,4] = 'LP1912'
$data1[297,0] = '15:22:17'
$data1[297,1] = '16:01'
$data1[297,2] = '10_OLMOS'
$data1[297,3] = 39
$data1[297,4] = 'LP1912'
$data1[298,0] = '15:53:28'
$data1[298,1] = '16:02'
$data1[298,2] = '16_SANTA ANA'
$data1[298,3] = 9
$data1[298,4] = 'LP1912'
$data1[299,0] = '15:53:28'
$data1[299,1] = '16:04'
$data1[299,2] = '23_HERNANDEZ'
$data1[299,3] = 11
$data1[299,4] = 'LP1912'
$data1[300,0] = '15:22:17'
$data1[300,1] = '16:05'
$data1[300,2] = '23_HERNANDEZ'
$data1[300,3] = 43
$data1[300,4] = 'LP1912'
$data1[301,0] = '14:56:20'
$data1[301,1] = '16:08'
$data1[301,2] = '14_ABASTO'
$data1[301,3] = 72
$data1[301,4] = 'LP1912'
$data1[302,0] = '14:45:17'
$data1[302,1] = '16:09'
$data1[302,2] = '14_ABASTO'
$data1[302,3] = 84
$data1[302,4] = 'LP1912'
$data1[303,0] = '14:24:16'
$data1[303,1] = '16:15'
$data1[303,2] = '225_C ROCA-H SUR'
$data1[303,3] = 111
$data1[303,4] = 'LP1912'
$data1[304,0] = '14:24:16'
$data1[304,1] = '16:20'
$data1[304,2] = '215C_EL PATO'
$data1[304,3] = 116
$data1[304,4] = 'LP1912'
$data1[305,0] = '14:24:16'
$data1[305,1] = '16:21'
$data1[305,2] = '26_HERNANDEZ'
$data1[305,3] = 117
$data1[305,4] = 'LP1912'
$data1[306,0] = '15:53:28'
$data1[306,1] = '16:29'
$data1[306,2] = '10_OLMOS'
$data1[306,3] = 36
$data1[306,4] = 'LP1912'
$data1[307,0] = '14:45:17'
$data1[307,1] = '16:30'
$data1[307,2] = '15_ABASTO'
$data1[307,3] = 105
$data1[307,4] = 'LP1912'
$data1[308,0] = '15:22:17'
$data1[308,1] = '16:32'
$data1[308,2] = '14_ABASTO'
$data1[308,3] = 70
$data1[308,4] = 'LP1912'
$data1[309,0] = '15:53:28'
$data1[309,1] = '16:34'
$data1[309,2] = '23_HERNANDEZ'
$data1[309,3] = 41
$data1[309,4] = 'LP1912'
$data1[310,0] = '15:53:28'
$data1[310,1] = '16:36'
$data1[310,2] = '11_ETCHEVERRY'
$data1[310,3] = 43
$data1[310,4] = 'LP1912'
$data1[311,0] = '15:22:17'
$data1[311,1] = '16:37'
$data1[311,2] = '11_ETCHEVERRY'
$data1[311,3] = 75
$data1[311,4] = 'LP1912'
$data1[312,0] = '15:22:17'
$data1[312,1] = '16:40'
$data1[312,2] = '17_ROMERO'
$data1[312,3] = 78
$data1[312,4] = 'LP1912'
$data1[313,0] = '14:56:20'
$data1[313,1] = '16:42'
$data1[313,2] = '16_P MOR-SANTA ANA'
$data1[313,3] = 106
$data1[313,4] = 'LP1912'
$data1[314,0] = '14:45:17'
$data1[314,1] = '16:43'
$data1[314,2] = '16_P MOR-SANTA ANA'
$data1[314,3] = 118
$data1[314,4] = 'LP1912'
$data1[315,0] = '14:45:17'
$data1[315,1] = '16:43'
$data1[315,2] = '225_GOMEZ'
$data1[315,3] = 118
$data1[315,4] = 'LP1912'
$data1[316,0] = '15:22:17'
$data1[316,1] = '16:48'
$data1[316,2] = '15_ABASTO'
$data1[316,3] = 86
$data1[316,4] = 'LP1912'
$data1[317,0] = '15:53:28'
$data1[317,1] = '16:50'
$data1[317,2] = '14_ABASTO'
$data1[317,3] = 57
$data1[317,4] = 'LP1912'
$data1[318,0] = '15:22:17'
$data1[318,1] = '16:56'
$data1[318,2] = '17_179 Y 38'
$data1[318,3] = 94
$data1[318,4] = 'LP1912'
$data1[319,0] = '15:22:17'
$data1[319,1] = '17:04'
$data1[319,2] = '215A_EL PATO'
$data1[319,3] = 102
$data1[319,4] = 'LP1912'
$data1[320,0] = '15:53:28'
$data1[320,1] = '17:21'
$data1[320,2] = '26_HERNANDEZ'
$data1[320,3] = 88
$data1[320,4] = 'LP1912'
$data1[321,0] = '15:53:28'
$data1[321,1] = '17:24'
$data1[321,2] = '84_COLONIA URQUIZA-ESC 49'
$data1[321,3] = 91
$data1[321,4] = 'LP1912'
$data1[322,0] = '15:53:28'
$data1[322,1] = '17:28'
$data1[322,2] = '14_ABASTO'
$data1[322,3] = 95
$data1[322,4] = 'LP1912'
$data1[323,0] = '15:53:28'
$data1[323,1] = '17:36'
$data1[323,2] = '27_EL RETIRO'
$data1[323,3] = 103
$data1[323,4] = 'LP1912'
$data1[324,0] = '15:53:28'
$data1[324,1] = '17:38'
$data1[324,2] = '17_ROMERO'
$data1[324,3] = 105
$data1[324,4] = 'LP1912'
$data1[325,0] = '15:53:28'
$data1[325,1] = '17:40'
$data1[325,2] = '215B_EL PATO'
$data1[325,3] = 107
$data1[325,4] = 'LP1912'
$data1[326,0] = '15:53:28'
$data1[326,1] = '17:50'
$data1[326,2] = '16_P MOR-167 Y 521'
$data1[326,3] = 117
$data1[326,4] = 'LP1912'
$data1[327,0] = '15:53:28'
$data1[327,1] = '17:52'
$data1[327,2] = '81_EL PELIGRO'
$data1[327,3] = 119
$data1[327,4] = 'LP1912'
$ws1.Range("A6:E333").Value = $data1

# --- Sheet: LP1912-215 (header) ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 15:53:28'
$ws2.Range("A3").Value = 'Total filas: 35'

# --- Sheet2 data block: A6:E40 (35 rows) ---
$data2 = New-Object 'object[,]' 35,5
$data2[0,0] = '03:45:25'
$data2[0,1] = '04:46'
$data2[0,2] = '215A_EL PATO'
$data2[0,3] = 61
$data2[0,4] = 'LP1912'
$data2[1,0] = '03:45:25'
$data2[1,1] = '05:34'
$data2[1,2] = '215B_EL PATO'
$data2[1,3] = 109
$data2[1,4] = 'LP1912'
$data2[2,0] = '04:18:02'
$data2[2,1] = '05:35'
$data2[2,2] = '215B_EL PATO'
$data2[2,3] = 77
$data2[2,4] = 'LP1912'
$data2[3,0] = '04:56:49'
$data2[3,1] = '06:11'
$data2[3,2] = '215A_EL PATO'
$data2[3,3] = 75
$data2[3,4] = 'LP1912'
$data2[4,0] = '04:18:02'
$data2[4,1] = '06:12'
$data2[4,2] = '215A_EL PATO'
$data2[4,3] = 114
$data2[4,4] = 'LP1912'
$data2[5,0] = '04:56:49'
$data2[5,1] = '06:46'
$data2[5,2] = '215C_EL PATO'
$data2[5,3] = 110
$data2[5,4] = 'LP1912'
$data2[6,0] = '05:26:08'
$data2[6,1] = '06:47'
$data2[6,2] = '215C_EL PATO'
$data2[6,3] = 81
$data2[6,4] = 'LP1912'
$data2[7,0] = '05:26:08'
$data2[7,1] = '07:11'
$data2[7,2] = '215A_EL PATO'
$data2[7,3] = 105
$data2[7,4] = 'LP1912'
$data2[8,0] = '06:55:02'
$data2[8,1] = '07:12'
$data2[8,2] = '215A_EL PATO'
$data2[8,3] = 17
$data2[8,4] = 'LP1912'
$data2[9,0] = '05:55:25'
$data2[9,1] = '07:51'
$data2[9,2] = '215D_EL PATO'
$data2[9,3] = 116
$data2[9,4] = 'LP1912'
$data2[10,0] = '06:55:02'
$data2[10,1] = '07:52'
$data2[10,2] = '215D_EL PATO'
$data2[10,3] = 57
$data2[10,4] = 'LP1912'
$data2[11,0] = '06:25:43'
$data2[11,1] = '08:23'
$data2[11,2] = '215B_EL PATO'
$data2[11,3] = 118
$data2[11,4] = 'LP1912'
$data2[12,0] = '07:19:29'
$data2[12,1] = '09:01'
$data2[12,2] = '215A_EL PATO'
$data2[12,3] = 102
$data2[12,4] = 'LP1912'
$data2[13,0] = '07:50:16'
$data2[13,1] = '09:02'
$data2[13,2] = '215A_EL PATO'
$data2[13,3] = 72
$data2[13,4] = 'LP1912'
$data2[14,0] = '07:50:16'
$data2[14,1] = '09:42'
$data2[14,2] = '215C_EL PATO'
$data2[14,3] = 112
$data2[14,4] = 'LP1912'
$data2[15,0] = '08:56:29'
$data2[15,1] = '10:26'
$data2[15,2] = '215A_EL PATO'
$data2[15,3] = 90
$data2[15,4] = 'LP1912'
$data2[16,0] = '08:32:09'
$data2[16,1] = '10:27'
$data2[16,2] = '215A_EL PATO'
$data2[16,3] = 115
$data2[16,4] = 'LP1912'
$data2[17,0] = '10:29:57'
$data2[17,1] = '11:01'
$data2[17,2] = '215C_EL PATO'
$data2[17,3] = 32
$data2[17,4] = 'LP1912'
$data2[18,0] = '09:35:26'
$data2[18,1] = '11:02'
$data2[18,2] = '215C_EL PATO'
$data2[18,3] = 87
$data2[18,4] = 'LP1912'
$data2[19,0] = '10:29:57'
$data2[19,1] = '11:50'
$data2[19,2] = '215B_EL PATO'
$data2[19,3] = 81
$data2[19,4] = 'LP1912'
$data2[20,0] = '10:59:49'
$data2[20,1] = '11:51'
$data2[20,2] = '215B_EL PATO'
$data2[20,3] = 52
$data2[20,4] = 'LP1912'
$data2[21,0] = '10:29:57'
$data2[21,1] = '12:20'
$data2[21,2] = '215A_EL PATO'
$data2[21,3] = 111
$data2[21,4] = 'LP1912'
$data2[22,0] = '12:21:08'
$data2[22,1] = '12:21'
$data2[22,2] = '215A_EL PATO'
$data2[22,3] = 0
$data2[22,4] = 'LP1912'
$data2[23,0] = '11:30:45'
$data2[23,1] = '13:13'
$data2[23,2] = '215D_EL PATO'
$data2[23,3] = 103
$data2[23,4] = 'LP1912'
$data2[24,0] = '12:21:08'
$data2[24,1] = '13:14'
$data2[24,2] = '215D_EL PATO'
$data2[24,3] = 53
$data2[24,4] = 'LP1912'
$data2[25,0] = '11:56:55'
$data2[25,1] = '13:50'
$data2[25,2] = '215A_EL PATO'
$data2[25,3] = 114
$data2[25,4] = 'LP1912'
$data2[26,0] = '12:21:08'
$data2[26,1] = '13:51'
$data2[26,2] = '215A_EL PATO'
$data2[26,3] = 90
$data2[26,4] = 'LP1912'
$data2[27,0] = '12:59:47'
$data2[27,1] = '14:19'
$data2[27,2] = '215C_EL PATO'
$data2[27,3] = 80
$data2[27,4] = 'LP1912'
$data2[28,0] = '12:21:08'
$data2[28,1] = '14:20'
$data2[28,2] = '215C_EL PATO'
$data2[28,3] = 119
$data2[28,4] = 'LP1912'
$data2[29,0] = '12:59:47'
$data2[29,1] = '14:58'
$data2[29,2] = '215B_EL PATO'
$data2[29,3] = 119
$data2[29,4] = 'LP1912'
$data2[30,0] = '14:45:17'
$data2[30,1] = '15:38'
$data2[30,2] = '215A_EL PATO'
$data2[30,3] = 53
$data2[30,4] = 'LP1912'
$data2[31,0] = '13:59:06'
$data2[31,1] = '15:39'
$data2[31,2] = '215A_EL PATO'
$data2[31,3] = 100
$data2[31,4] = 'LP1912'
$data2[32,0] = '14:24:16'
$data2[32,1] = '16:20'
$data2[32,2] = '215C_EL PATO'
$data2[32,3] = 116
$data2[32,4] = 'LP1912'
$data2[33,0] = '15:22:17'
$data2[33,1] = '17:04'
$data2[33,2] = '215A_EL PATO'
$data2[33,3] = 102
$data2[33,4] = 'LP1912'
$data2[34,0] = '15:53:28'
$data2[34,1] = '17:40'
$data2[34,2] = '215B_EL PATO'
$data2[34,3] = 107
$data2[34,4] = 'LP1912'
$ws2.Range("A6:E40").Value = $data2

# --- Sheet: 6203-6173 (header only) ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 15:53:28'

Write-Output "done"